# Weekly update: insert a new data row for Coliflor (Macroferia Regional de Talca)
# above the current row 50, shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 50; this shifts rows 50..112 down to 51..113
# and copies formatting (incl. the date number format on column D) from row 50.
$ws.Rows.Item(50).Insert()

# Fill in the new row 50 with this week's record.
$ws.Cells.Item(50, 1).Value = 5
$ws.Cells.Item(50, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(50, 3).Value = "Maule"
$ws.Cells.Item(50, 4).Value = 44413
$ws.Cells.Item(50, 5).Value = 7
$ws.Cells.Item(50, 6).Value = 100112008
$ws.Cells.Item(50, 7).Value = "Coliflor"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Segunda"
$ws.Cells.Item(50, 10).Value = 3000
$ws.Cells.Item(50, 11).Value = 500
$ws.Cells.Item(50, 12).Value = 500
$ws.Cells.Item(50, 13).Value = 500
$ws.Cells.Item(50, 14).Value = '$/unidad'
$ws.Cells.Item(50, 15).Value = "Región del Maule"
$ws.Cells.Item(50, 16).Value = 500
$ws.Cells.Item(50, 17).Value = 1
$ws.Cells.Item(50, 18).Value = "Hortaliza"
